$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Output ("Sheet name: " + $ws.Name)
